$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.507.12'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '1.819.17'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'315.47"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.87%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = "'0.5065"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.23%  '
$ws.Range('D8').Value = "'0.3840"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.49%  '
$ws.Range('D9').Value = "'0.08479"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.81%  '
$ws.Range('D10').Value = "'41.93"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'1.107"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('D12').Value = "'6.407"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').Value = "'21.04"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').Value = "'7.487"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = '1.813.39'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').Value = "'0.00001143"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.89%  '
$ws.Range('D18').Value = "'93.24"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').Value = "'0.06702"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').Value = "'17.70"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').Value = "'1.002"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = "'6.070"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').Value = '28.526.48'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').Value = "'11.41"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('D26').Value = "'21.24"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('D27').Value = "'158.87"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('D28').Value = '2.022.63'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').Value = "'2.384"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').Value = "'126.04"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('D31').Value = "'1.104"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('D32').Value = "'0.1076"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.28%  '
$ws.Range('D33').Value = "'5.743"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = "'3.689"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = "'0.07358"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('D36').Value = "'0.2224"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('D37').Value = "'0.02360"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = "'5.201"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('D39').Value = "'8.715"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.29%  '
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').Value = "'11.25"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.55%  '
$ws.Range('D42').Value = "'1.192"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').Value = "'1.403"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('D44').Value = "'13.59"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').Value = "'3.750"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.74%  '
$ws.Range('D46').Value = "'0.5902"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').Value = "'125.45"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').Value = "'1.988"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('D49').Value = "'1.195"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').Value = "'0.06993"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').Value = "'74.06"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.48%  '
